$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column for existing rows 2-499: 45192 -> 45202
$ws.Range("C2:C499").Value = 45202

# 2. Row 499 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(499).RowHeight = 15

# 3. Append new row 500 with the new avverkningsanmälan entry
$ws.Range("A500").Value = "A 46165-2023"

$ws.Range("B500").Value = 45196
$ws.Range("B500").NumberFormat = "YYYY-MM-DD"

$ws.Range("C500").Value = 45202
$ws.Range("C500").NumberFormat = "YYYY-MM-DD"

$ws.Range("D500").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E500").Value = "GISLAVED"

$ws.Range("G500").Value = 4
$ws.Range("H500").Value = 0
$ws.Range("I500").Value = 0
$ws.Range("J500").Value = 0
$ws.Range("K500").Value = 0
$ws.Range("L500").Value = 0
$ws.Range("M500").Value = 0
$ws.Range("N500").Value = 0
$ws.Range("O500").Value = 0
$ws.Range("P500").Value = 0
$ws.Range("Q500").Value = 0

$ws.Range("R500").Value = ""
$ws.Range("R500").WrapText = $true
